# printLabelSaoHanTemplate.docx — "improve print labels process, print 2 labels/person"
#
# Changes applied:
#  1. Bump the (empty) paragraph's mark run-properties font size 14pt -> 15.5pt
#     (w:sz/w:szCs 28 -> 31 half-points) so two labels fit per printed sheet.
#  2. Resize the page / margins for the new 2-up label stock
#     (pgSz 11338x23811 -> 11282x20160 twips; top margin 170 -> 227,
#      left margin 400 -> 386 twips).
#  3. Mark the "Default Paragraph Font" style as a quick/recommended style.

$d = $word.ActiveDocument

# --- 1. Paragraph-mark font size -------------------------------------------
# The document body is a single, completely empty paragraph — there is no
# run for the COM Font setters to attach to, so Word (and this host) will
# silently ignore a direct `Paragraphs(1).Range.Font.Size = ...` on it.
# Route around that the same way a user typing in Word would: type a
# placeholder character, size the range that includes the paragraph mark
# (which folds the size into the paragraph mark's run properties), then
# delete the placeholder again so the paragraph is left empty exactly as
# before, only with the new mark formatting retained.
$sel = $word.Selection
$sel.TypeText("X")

$markRange = $d.Paragraphs.Item(1).Range
$markRange.Font.Size = 15.5
$markRange.Font.SizeBi = 15.5

$placeholder = $d.Range(0, 1)
$placeholder.Delete()

# --- 2. Page size / margins --------------------------------------------------
$ps = $d.PageSetup
$ps.PageWidth = 564.1      # 11282 twips
$ps.PageHeight = 1008.0    # 20160 twips
$ps.TopMargin = 11.35      # 227 twips
$ps.LeftMargin = 19.3      # 386 twips
# Right/bottom margins (538 / 1440 twips) are unchanged by the diff.

# --- 3. "Default Paragraph Font" becomes a quick style -----------------------
$dpf = $d.Styles.Item("Default Paragraph Font")
$dpf.QuickStyle = $true
